# Update Mappings 22 Ontologies
# - Swap the CHEBI_51086 / CHEBI_24432 rows (rows 6 and 7) back into sync
#   with the refreshed ontology mapping order.
# - Add a new "ChEBI_DEF" column (F) holding the ChEBI definition payload
#   for each mapped row (only the first data row has a real definition;
#   the rest default to an empty list "[]").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-sync rows 6 and 7 (CHEBI_51086 <-> CHEBI_24432) ---------------
$ws.Range("B6").Value = "http://purl.obolibrary.org/obo/CHEBI_24432"
$ws.Range("C6").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_24432'}"
$ws.Range("D6").Value = "http://purl.obolibrary.org/obo/CHEBI_24432"
$ws.Range("E6").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_24432'}"

$ws.Range("B7").Value = "http://purl.obolibrary.org/obo/CHEBI_51086"
$ws.Range("C7").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_51086'}"
$ws.Range("D7").Value = "http://purl.obolibrary.org/obo/CHEBI_51086"
$ws.Range("E7").Value = "{'iri': 'http://purl.obolibrary.org/obo/CHEBI_51086'}"

# --- 2. Add the new ChEBI_DEF column (F) ---------------------------------
# Copy the header formatting from E1 (bold/centered header style) onto F1.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "ChEBI_DEF"

$iupacDef = "['Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]']"

$ws.Range("F2").Value = $iupacDef

for ($r = 3; $r -le 46; $r++) {
    $ws.Cells.Item($r, 6).Value = "[]"
}
